$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44748
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 2300
$ws.Range("O2").Value = 2300
$ws.Range("P2").Value = 2300
$ws.Range("S2").Value = 2300

$ws.Range("D3").Value = 44749
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 2300
$ws.Range("O3").Value = 2300
$ws.Range("P3").Value = 2300
$ws.Range("S3").Value = 2300

$ws.Range("D4").Value = 45075
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 3200
$ws.Range("O4").Value = 3200
$ws.Range("P4").Value = 3200
$ws.Range("S4").Value = 3200

$ws.Range("D5").Value = 44763
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 2300
$ws.Range("O5").Value = 2300
$ws.Range("P5").Value = 2300
$ws.Range("S5").Value = 2300

$ws.Range("D6").Value = 44762
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 2300
$ws.Range("O6").Value = 2300
$ws.Range("P6").Value = 2300
$ws.Range("S6").Value = 2300

$ws.Range("D7").Value = 44753
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 2300
$ws.Range("O7").Value = 2300
$ws.Range("P7").Value = 2300
$ws.Range("S7").Value = 2300

$ws.Range("D8").Value = 45086
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 2600
$ws.Range("O8").Value = 2600
$ws.Range("P8").Value = 2600
$ws.Range("S8").Value = 2600

$ws.Range("D9").Value = 44811
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 2500
$ws.Range("O9").Value = 2500
$ws.Range("P9").Value = 2500
$ws.Range("S9").Value = 2500

$ws.Range("D10").Value = 45062
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 3200
$ws.Range("O10").Value = 3200
$ws.Range("P10").Value = 3200
$ws.Range("S10").Value = 3200

$ws.Range("D11").Value = 45068
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 3250
$ws.Range("O11").Value = 3250
$ws.Range("P11").Value = 3250
$ws.Range("S11").Value = 3250

$ws.Range("D12").Value = 44760
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 2300
$ws.Range("O12").Value = 2300
$ws.Range("P12").Value = 2300
$ws.Range("S12").Value = 2300

$ws.Range("D13").Value = 44435
$ws.Range("M13").Value = 130
$ws.Range("N13").Value = 1300
$ws.Range("O13").Value = 1300
$ws.Range("P13").Value = 1300
$ws.Range("S13").Value = 1300

$ws.Range("D14").Value = 45054
$ws.Range("M14").Value = 25
$ws.Range("N14").Value = 2500
$ws.Range("O14").Value = 2500
$ws.Range("P14").Value = 2500
$ws.Range("S14").Value = 2500

$ws.Range("D15").Value = 44417
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 1200
$ws.Range("O15").Value = 1200
$ws.Range("P15").Value = 1200
$ws.Range("S15").Value = 1200

$ws.Range("D16").Value = 44476
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 1200
$ws.Range("O16").Value = 1200
$ws.Range("P16").Value = 1200
$ws.Range("S16").Value = 1200

$ws.Range("D17").Value = 45076
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 2600
$ws.Range("O17").Value = 2600
$ws.Range("P17").Value = 2600
$ws.Range("S17").Value = 2600

$ws.Range("D18").Value = 45044
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 3500
$ws.Range("O18").Value = 3500
$ws.Range("P18").Value = 3500
$ws.Range("S18").Value = 3500

$ws.Range("D19").Value = 44812
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 2500
$ws.Range("O19").Value = 2500
$ws.Range("P19").Value = 2500
$ws.Range("S19").Value = 2500

$ws.Range("D20").Value = 44357
$ws.Range("M20").Value = 35
$ws.Range("N20").Value = 1000
$ws.Range("O20").Value = 1000
$ws.Range("P20").Value = 1000
$ws.Range("S20").Value = 1000

$ws.Range("D21").Value = 44438
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 1200
$ws.Range("O21").Value = 1200
$ws.Range("P21").Value = 1200
$ws.Range("S21").Value = 1200

$ws.Range("D22").Value = 44343
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 1300
$ws.Range("O22").Value = 1300
$ws.Range("P22").Value = 1300
$ws.Range("S22").Value = 1300

$ws.Range("D23").Value = 45079
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 2600
$ws.Range("O23").Value = 2600
$ws.Range("P23").Value = 2600
$ws.Range("S23").Value = 2600

$ws.Range("D24").Value = 45042
$ws.Range("M24").Value = 25
$ws.Range("N24").Value = 3500
$ws.Range("O24").Value = 3500
$ws.Range("P24").Value = 3500
$ws.Range("S24").Value = 3500

$ws.Range("D25").Value = 44431
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 1300
$ws.Range("O25").Value = 1300
$ws.Range("P25").Value = 1300
$ws.Range("S25").Value = 1300

$ws.Range("D26").Value = 44418
$ws.Range("M26").Value = 40
$ws.Range("N26").Value = 1200
$ws.Range("O26").Value = 1200
$ws.Range("P26").Value = 1200
$ws.Range("S26").Value = 1200

$ws.Range("D27").Value = 45090
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 2600
$ws.Range("O27").Value = 2600
$ws.Range("P27").Value = 2600
$ws.Range("S27").Value = 2600

$ws.Range("D28").Value = 45041
$ws.Range("M28").Value = 80
$ws.Range("N28").Value = 3500
$ws.Range("O28").Value = 3500
$ws.Range("P28").Value = 3500
$ws.Range("S28").Value = 3500

$ws.Range("D29").Value = 45085
$ws.Range("M29").Value = 40
$ws.Range("N29").Value = 2600
$ws.Range("O29").Value = 2600
$ws.Range("P29").Value = 2600
$ws.Range("S29").Value = 2600

$ws.Range("D30").Value = 44405
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 1200
$ws.Range("O30").Value = 1200
$ws.Range("P30").Value = 1200
$ws.Range("S30").Value = 1200

$ws.Range("D31").Value = 44432
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 1300
$ws.Range("O31").Value = 1300
$ws.Range("P31").Value = 1300
$ws.Range("S31").Value = 1300

$ws.Range("D32").Value = 45055
$ws.Range("M32").Value = 25
$ws.Range("N32").Value = 2800
$ws.Range("O32").Value = 2800
$ws.Range("P32").Value = 2800
$ws.Range("S32").Value = 2800

$ws.Range("D33").Value = 44424
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 1200
$ws.Range("O33").Value = 1200
$ws.Range("P33").Value = 1200
$ws.Range("S33").Value = 1200

$ws.Range("D34").Value = 44473
$ws.Range("M34").Value = 120
$ws.Range("N34").Value = 1200
$ws.Range("O34").Value = 1200
$ws.Range("P34").Value = 1200
$ws.Range("S34").Value = 1200

$ws.Range("D35").Value = 44830
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 2500
$ws.Range("O35").Value = 2500
$ws.Range("P35").Value = 2500
$ws.Range("S35").Value = 2500
